$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Samples")
$ws.Range("G2").Value = "generic"
$ws.Range("G3").Select()
